$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the missing Friday hours for the week commencing 12/03/2018 (row 10)
$ws.Range("B10").Value = 4.25

# Update the active cell selection to reflect the last-edited cell
$ws.Range("F15").Select()
